$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# 1) Status text: "Ready for handoff" -> "In Translation"
#    Occurs on the Overview sheet (summary columns for each locale) and
#    on each per-locale sheet's "Status" column.
# ----------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# ----------------------------------------------------------------------
# 2) Narrow the "Status" related columns.
#    Overview!E:F and the "Status" column (C) on each locale sheet shrink
#    from ~17.22 chars to ~13.41 chars.
# ----------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
